$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Singlecard")

# Update row 2 values
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "01-08-1988"
$ws.Range("D2").Value = 33
$ws.Range("F2").Value = "МБОУ средняя общеобразовательная школа с.Виноградное"
$ws.Range("H2").Value = "Экономист"
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 10
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "2019"
$ws.Range("N2").Value = "vinogradovskayasos1h@mail.ru"
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "89991118888"

# Remove the third row entirely (Саидов Саид Саидович ...)
$ws.Rows.Item(3).Delete()
